# Revamp PricingDownload: append the new "Equal Exchange - One World" line
# item as row 4 of the pricing sheet (SKU, Item, Quantity, Cost Per, Total
# Cost), matching the existing sheet's convention of storing every value as
# text (not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row's cells to Text format first so values like "10400" /
# "3" / "71.50" / "214.50" are stored as strings (consistent with the rest
# of the sheet) instead of being auto-coerced into numbers.
$ws.Range("A4:E4").NumberFormat = "@"

$ws.Range("A4").Value = "10400"
$ws.Range("B4").Value = "Equal Exchange - One World"
$ws.Range("C4").Value = "3"
$ws.Range("D4").Value = "71.50"
$ws.Range("E4").Value = "214.50"

# Drop back to the default "Normal" style so no stray text-format styling
# is left attached to the new cells.
$ws.Range("A4:E4").Style = "Normal"
